# --- 对外直接投资流量.xlsx update -------------------------------------
# Commit intent (from the OOXML diff): the table of "对外直接投资流量"
# (outbound FDI flow) by industry drops the earliest three years (2007-
# 2009), replaces the legacy high-precision 2020 figures with the
# finalized/rounded ones, and appends a new 2021 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the 2007/2008/2009 rows (rows 2-4). Everything below shifts up
#    by three rows, so the former 2010 row becomes row 2, ..., and the
#    former 2020 row (was row 15) becomes row 12.
$ws.Rows("2:4").Delete()

# 2) The former row 15 (2020年) carried long decimal figures; replace
#    them with the finalized, rounded values now that it lives at row 12.
$ws.Range("B12").Value2 = 623320
$ws.Range("C12").Value2 = 11841
$ws.Range("D12").Value2 = 918718
$ws.Range("F12").Value2 = 107864
$ws.Range("G12").Value2 = 2583821
$ws.Range("H12").Value2 = 63767
$ws.Range("I12").Value2 = 15371026
$ws.Range("J12").Value2 = 216078
$ws.Range("K12").Value2 = 809455
$ws.Range("L12").Value2 = 518603
$ws.Range("M12").Value2 = 2299764
$ws.Range("N12").Value2 = 13004
$ws.Range("O12").Value2 = -213383
$ws.Range("P12").Value2 = 15671
$ws.Range("Q12").Value2 = 577031
$ws.Range("R12").Value2 = 373465
$ws.Range("S12").Value2 = 3872562
$ws.Range("T12").Value2 = 613126
$ws.Range("U12").Value2 = 1966318

# 3) Append the new 2021年 row (row 13). Clone row 12's formatting for
#    the year-label cell (bold, bordered, centered) before filling values.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value2 = "2021年"
$ws.Range("B13").Value2 = 1222621
$ws.Range("C13").Value2 = 26933
$ws.Range("D13").Value2 = 513591
$ws.Range("F13").Value2 = 93075
$ws.Range("G13").Value2 = 2686673
$ws.Range("H13").Value2 = 33877
$ws.Range("I13").Value2 = 17881932
$ws.Range("J13").Value2 = 180948
$ws.Range("K13").Value2 = 461908
$ws.Range("L13").Value2 = 409785
$ws.Range("M13").Value2 = 2815201
$ws.Range("N13").Value2 = 2825
$ws.Range("O13").Value2 = 8773
$ws.Range("P13").Value2 = 22494
$ws.Range("Q13").Value2 = 438908
$ws.Range("R13").Value2 = 507213
$ws.Range("S13").Value2 = 4935732
$ws.Range("T13").Value2 = 841498
$ws.Range("U13").Value2 = 2679879
